# Updates the cryptocurrency price/volume table with freshly scraped values.
# Commit: "Updated cryptos list on Tue Nov  7 18:16:20 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cell references (row-major, matching the sheet layout).
$refs = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "D22",
    "E22",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D28",
    "E28",
    "D29",
    "E29",
    "E30",
    "D31",
    "E31",
    "D32",
    "E32",
    "D33",
    "E33",
    "D34",
    "E34",
    "E35",
    "D36",
    "E36",
    "D37",
    "E37",
    "D38",
    "E38",
    "B39",
    "C39",
    "D39",
    "E39",
    "B40",
    "C40",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "B45",
    "C45",
    "D45",
    "E45",
    "B46",
    "C46",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)

# New text for each cell above. Values that would otherwise be auto-parsed
# as numbers are apostrophe-prefixed so Excel keeps them as literal text
# (matching the source data, e.g. "34.872.12" or "0.998" as plain strings).
$vals = @(
    "35.282.87",
    "  +0.07%  ",
    "1.886.10",
    "  -1.31%  ",
    "  -1.11%  ",
    "'246.45",
    "  -2.72%  ",
    "'0.686",
    "  -4.66%  ",
    "'0.998",
    "  -1.20%  ",
    "'43.60",
    "  +7.52%  ",
    "'0.348",
    "  -3.24%  ",
    "'0.0739",
    "  -4.02%  ",
    "'0.0970",
    "  -1.89%  ",
    "'13.08",
    "  +3.14%  ",
    "2.155.46",
    "  -1.64%  ",
    "'0.719",
    "  +0.23%  ",
    "'4.91",
    "  +0.11%  ",
    "1.879.63",
    "  -3.67%  ",
    "35.150.45",
    "  -0.44%  ",
    "'73.02",
    "  -1.52%  ",
    "0.0₃0819",
    "  -3.84%  ",
    "'244.58",
    "  +0.50%  ",
    "'12.78",
    "  -1.48%  ",
    "'4.96",
    "  -2.21%  ",
    "  -0.47%  ",
    "'2.52",
    "  +6.96%  ",
    "'2.18",
    "  -9.71%  ",
    "'164.13",
    "  -1.77%  ",
    "'8.47",
    "  -1.54%  ",
    "'18.22",
    "  -2.48%  ",
    "'0.126",
    "  -4.26%  ",
    "  -0.01%  ",
    "'1.78",
    "  +10.14%  ",
    "'4.24",
    "  -2.00%  ",
    "'0.0578",
    "  -0.79%  ",
    "'4.21",
    "  +0.21%  ",
    "  -1.10%  ",
    "'0.839",
    "  -8.37%  ",
    "'1.98",
    "  -2.23%  ",
    "'1.48",
    "  -25.95%  ",
    "InjectiveProtocol",
    "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj",
    "'17.31",
    "  +0.90%  ",
    "Aave",
    "https://coinranking.com/coin/ixgUfzmLR+aave-aave",
    "'98.44",
    "  +1.70%  ",
    "'0.0667",
    "  +3.94%  ",
    "'0.0212",
    "  -2.64%  ",
    "'1.09",
    "  -2.68%  ",
    "1.288.95",
    "  -3.62%  ",
    "RenderToken",
    "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr",
    "'2.35",
    "  -3.33%  ",
    "Cronos",
    "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro",
    "'0.0817",
    "  +10.41%  ",
    "'2.40",
    "  -1.22%  ",
    "'2.74",
    "  -1.04%  ",
    "'12.04",
    "  -0.36%  ",
    "'6.37",
    "  -6.16%  ",
    "'43.08",
    "  -4.66%  "
)

for ($i = 0; $i -lt $refs.Count; $i++) {
    $cell = $ws.Range($refs[$i])
    $v = $vals[$i]
    $cell.Value = $v
    if ($v.StartsWith("'")) {
        # Drop the quote-prefix style Excel applied for the text coercion
        # above so only the cell content (not its style) changes.
        $cell.Style = "Normal"
    }
}

